# Update "想去人数" (column F) values across all sheets to match regenerated source data
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 751
$ws.Range("F3").Value = 14154
$ws.Range("F4").Value = 14155
$ws.Range("F5").Value = 14225
$ws.Range("F7").Value = 1394
$ws.Range("F8").Value = 5866
$ws.Range("F14").Value = 1535
$ws.Range("F15").Value = 434
$ws.Range("F16").Value = 2130
$ws.Range("F17").Value = 1190
$ws.Range("F18").Value = 1814
$ws.Range("F19").Value = 913
$ws.Range("F21").Value = 2269
$ws.Range("F22").Value = 560
$ws.Range("F23").Value = 806
$ws.Range("F24").Value = 3309
$ws.Range("F26").Value = 309
$ws.Range("F27").Value = 2379
$ws.Range("F28").Value = 584
$ws.Range("F31").Value = 1778
$ws.Range("F32").Value = 1067
$ws.Range("F33").Value = 1375
$ws.Range("F34").Value = 100
$ws.Range("F35").Value = 144
$ws.Range("F36").Value = 4764
$ws.Range("F37").Value = 4812
$ws.Range("F38").Value = 299
$ws.Range("F40").Value = 670
$ws.Range("F42").Value = 3282
$ws.Range("F43").Value = 42
$ws.Range("F44").Value = 921
$ws.Range("F46").Value = 100
$ws.Range("F47").Value = 71
$ws.Range("F48").Value = 4417
$ws.Range("F49").Value = 571
$ws.Range("F50").Value = 287

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 119
$ws.Range("F19").Value = 89
$ws.Range("F22").Value = 54

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 7523
$ws.Range("F3").Value = 233
$ws.Range("F4").Value = 743

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 7523
$ws.Range("F3").Value = 751
$ws.Range("F4").Value = 233
$ws.Range("F5").Value = 743
$ws.Range("F6").Value = 14155
$ws.Range("F7").Value = 14225
$ws.Range("F9").Value = 1394
$ws.Range("F10").Value = 5866
$ws.Range("F12").Value = 119
$ws.Range("F15").Value = 1535
$ws.Range("F16").Value = 434
$ws.Range("F17").Value = 1190
$ws.Range("F18").Value = 1814
$ws.Range("F19").Value = 913
$ws.Range("F21").Value = 3309
$ws.Range("F22").Value = 309
$ws.Range("F23").Value = 2379
$ws.Range("F24").Value = 584
$ws.Range("F27").Value = 1778
$ws.Range("F31").Value = 1067
$ws.Range("F32").Value = 1375
$ws.Range("F33").Value = 100
$ws.Range("F34").Value = 4764
$ws.Range("F35").Value = 4812
$ws.Range("F36").Value = 299
$ws.Range("F38").Value = 670
$ws.Range("F40").Value = 3282
$ws.Range("F41").Value = 921
$ws.Range("F43").Value = 100
$ws.Range("F45").Value = 71
$ws.Range("F46").Value = 4417
$ws.Range("F47").Value = 571
$ws.Range("F48").Value = 287
